$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "PARQUE DE GORONGOSA"
$ws.Range("E2").Value = "TECNICO DE POUPANCA"
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = " 11/9/2022"
$ws.Range("H2").Value = 2
$ws.Range("K2").Value = ""

# Row 3
$ws.Range("D3").Value = "PARQUE DE GORONGOSA"
$ws.Range("E3").Value = "TECNICO DE CAFE"
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = " 11/9/2022"
$ws.Range("H3").Value = 2
$ws.Range("K3").Value = ""

# Row 4
$ws.Range("C4").Value = "FEMININO"
$ws.Range("D4").Value = "PARQUE DE GORONGOSA"
$ws.Range("E4").Value = "SUPERVISORA DE ARTESANATO"
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = "  9/8/2022"
$ws.Range("H4").Value = 1

# Row 6
$ws.Range("D6").Value = "PARQUE DE GORONGOSA"
$ws.Range("E6").Value = "TECNICO DE CAFE"
$ws.Range("F6").Value = ""
$ws.Range("G6").Value = "  9/8/2022"
$ws.Range("H6").Value = 1
$ws.Range("K6").Value = ""
